{"js": "// Apply the \"Final edits for publication\" changes to the Zika manuscript.\n//\n// 1. Title-case several words in the title (first paragraph):\n//      impact -> Effect, discontinuing -> Discontinuing, universal -> Universal,\n//      screening -> Screening, donated -> Donated, blood -> Blood,\n//      virus -> Virus, states -> States\n// 2. \"1483 years\" -> \"1484 years\" in the congenital-Zika-syndrome sentence.\n// 3. Rewording of the outbreak-resurgence sentence in the Discussion.\n// 4. Rewording of the \"Enhanced collaboration...\" closing sentence.\n\nconst body = context.document.body;\n\n// --- 1. Title paragraph word replacements -----------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\n\nconst titleWordReplacements = [\n  [\"impact\", \"Effect\"],\n  [\"discontinuing\", \"Discontinuing\"],\n  [\"universal\", \"Universal\"],\n  [\"screening\", \"Screening\"],\n  [\"donated\", \"Donated\"],\n  [\"blood\", \"Blood\"],\n  [\"virus\", \"Virus\"],\n  [\"states\", \"States\"],\n];\n\nfor (const [oldWord, newWord] of titleWordReplacements) {\n  const found = titlePara.search(oldWord, { matchCase: true, matchWholeWord: true });\n  found.load(\"items\");\n  await context.sync();\n  for (const r of found.items) {\n    r.insertText(newWord, \"Replace\");\n  }\n  await context.sync();\n}\n\n// --- 2. Numeric correction: 1483 -> 1484 --------------------------------\nconst yearsHit = body.search(\"1483\", { matchCase: true });\nyearsHit.load(\"items\");\nawait context.sync();\nfor (const r of yearsHit.items) {\n  r.insertText(\"1484\", \"Replace\");\n}\nawait context.sync();\n\n// --- 3. Resurgence sentence rewording -----------------------------------\nconst resurgenceOld =\n  \"While Zika resurgence is possible, low endemic levels of Zika and established herd immunity in many regions make a larger outbreak unlikely.\";\nconst resurgenceNew =\n  \"While resurgence is possible, widespread immune experience and low endemic levels of Zika following the earlier pandemic make a larger outbreak unlikely.\";\n\nconst resurgenceHit = body.search(resurgenceOld, { matchCase: true });\nresurgenceHit.load(\"items\");\nawait context.sync();\nfor (const r of resurgenceHit.items) {\n  r.insertText(resurgenceNew, \"Replace\");\n}\nawait context.sync();\n\n// --- 4. Closing-paragraph collaboration sentence rewording ---------------\nconst collabOld =\n  \"Enhanced collaboration could improve understanding between policymakers and modelers, increasing the policy relevance and utilization of model-based reports.\";\nconst collabNew =\n  \"Collaboration could improve understanding between policymakers and modelers the policy relevance and utilization of model-based reports.\";\n\nconst collabHit = body.search(collabOld, { matchCase: true });\ncollabHit.load(\"items\");\nawait context.sync();\nfor (const r of collabHit.items) {\n  r.insertText(collabNew, \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Apply the \"Final edits for publication\" changes to the Zika manuscript.\n#\n# 1. Title-case several words in the title (paragraph 1):\n#      impact -> Effect, discontinuing -> Discontinuing, universal -> Universal,\n#      screening -> Screening, donated -> Donated, blood -> Blood,\n#      virus -> Virus, states -> States\n# 2. \"1483 years\" -> \"1484 years\" in the congenital-Zika-syndrome sentence.\n# 3. Rewording of the outbreak-resurgence sentence in the Discussion.\n# 4. Rewording of the \"Enhanced collaboration...\" closing sentence.\n\n$d = $word.ActiveDocument\n\n# --- 1. Title paragraph word replacements -----------------------------\n$titleRange = $d.Paragraphs(1).Range\n\n$titleWordReplacements = @(\n  @(\"impact\", \"Effect\"),\n  @(\"discontinuing\", \"Discontinuing\"),\n  @(\"universal\", \"Universal\"),\n  @(\"screening\", \"Screening\"),\n  @(\"donated\", \"Donated\"),\n  @(\"blood\", \"Blood\"),\n  @(\"virus\", \"Virus\"),\n  @(\"states\", \"States\")\n)\n\nforeach ($pair in $titleWordReplacements) {\n  $oldWord = $pair[0]\n  $newWord = $pair[1]\n  # MatchCase:=True, MatchWholeWord:=True, Replace:=wdReplaceOne(1) -\n  # scoped to the title paragraph's own Range so none of the many later\n  # lower-case occurrences of these common words are touched.\n  $titleRange.Find.Execute($oldWord, $true, $true, $false, $false, $false, $true, 1, $false, $newWord, 1) | Out-Null\n}\n\n# --- 2. Numeric correction: 1483 -> 1484 --------------------------------\n$d.Content.Find.Execute(\"1483\", $true, $false, $false, $false, $false, $true, 1, $false, \"1484\", 1) | Out-Null\n\n# --- 3. Resurgence sentence rewording -----------------------------------\n$resurgenceOld = \"While Zika resurgence is possible, low endemic levels of Zika and established herd immunity in many regions make a larger outbreak unlikely.\"\n$resurgenceNew = \"While resurgence is possible, widespread immune experience and low endemic levels of Zika following the earlier pandemic make a larger outbreak unlikely.\"\n$d.Content.Find.Execute($resurgenceOld, $true, $false, $false, $false, $false, $true, 1, $false, $resurgenceNew, 1) | Out-Null\n\n# --- 4. Closing-paragraph collaboration sentence rewording ---------------\n$collabOld = \"Enhanced collaboration could improve understanding between policymakers and modelers, increasing the policy relevance and utilization of model-based reports.\"\n$collabNew = \"Collaboration could improve understanding between policymakers and modelers the policy relevance and utilization of model-based reports.\"\n$d.Content.Find.Execute($collabOld, $true, $false, $false, $false, $false, $true, 1, $false, $collabNew, 1) | Out-Null\n"}
